$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "load time" formulas for the
# "Pubmed - Mesh + Supp." / SemmedDB table (rows 14-18, column D).
$ws.Range("D14").Formula = "=(203.385+203.583+204.463+204.39+203.633)/5"
$ws.Range("D15").Formula = "=(6.11628+6.10135+6.11597+6.1265+6.14355)/5"
$ws.Range("D16").Formula = "=(488.498+482.918+485.014+484.759+481.954)/5"
$ws.Range("D17").Formula = "=(2070.61+2136.17+2025.84+1536.91+1697.72)/5"
$ws.Range("D18").Formula = "=(159.598+162.842+160.258+163.46+164.336)/5"

# D19 already holds SUM(D14:D18); it will recalculate automatically.

# Move the active selection from E19 to E13.
$ws.Range("E13").Select()

# Slightly widen column A (and the rest of the default columns, which
# share the same column definition) from ~11.52 to ~11.96 characters.
$ws.Columns.Item(1).ColumnWidth = 11.17
